# "Running failed search test cases"
# - Remove the now-obsolete TestCase_B15 / B16 / B17 rows (10 MORE button tests)
# - Update a few test-case descriptions
# - Mark TestCase_B8 as PASS (it was re-run)
# - Reset Runmode to N for every test that isn't being (re-)run this pass,
#   keeping Y only for the search test cases that were actually run
#   (B6, B8, B11, B18, B19)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# TestCase_B15, B16 and B17 ("10 MORE button...") used to live on rows 16-18;
# they've been removed from the suite entirely.
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(16).Delete()

# TestCase_B8 (row 9): description rewritten, now verified as passing.
$ws.Range("B9").Value = "To verify that number of displayed documents gets increased as and when user scrolls down the search results page"
$ws.Range("D9").Value = "PASS"

# TestCase_B11 (row 12): dropped the "search," qualifier from the description.
$ws.Range("B12").Value = "To verify that sorting and filtering are retained when user navigates back to search results page from record view page"

# TestCase_B18 (row 16, after the B15-B17 rows were removed): dropped "Views,".
$ws.Range("B16").Value = "To verify that Times cited and Comments fields are getting displayed for each document in search results page"

# Runmode column: only the test cases actually exercised in this run stay "Y";
# everything else is reset to "N".
$runRows = @(7, 9, 12, 16, 17)
for ($r = 2; $r -le 25; $r++) {
    if ($runRows -contains $r) {
        $ws.Cells.Item($r, 3).Value = "Y"
    } else {
        $ws.Cells.Item($r, 3).Value = "N"
    }
}

# Selection marker left on C18 (as captured in the saved workbook).
$ws.Range("C18").Select()
